$d = $word.ActiveDocument

# 1. Replace phone number
$d.Content.Find.Execute("+91-894-064-9404", $true, $false, $false, $false, $false,
                         $true, 1, $false, "+32-494-80-87-16", 2)

Write-Host "Done"
